$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1411.0952
$ws.Range("I28").Value = 1013.7059
$ws.Range("K28").Value = 1013.7059
$ws.Range("M28").Value = -528.7059

$ws.Range("H32").Value = 14104.4
$ws.Range("I32").Value = 12880
$ws.Range("J32").Value = 19002
$ws.Range("K32").Value = 12880
$ws.Range("L32").Value = 19002
$ws.Range("M32").Value = -12554
$ws.Range("N32").Value = -19654

$ws.Range("H40").Value = 6074.75
$ws.Range("J40").Value = 8599.4
$ws.Range("L40").Value = 8599.4
$ws.Range("N40").Value = -8949.4

$ws.Range("H43").Value = 5063
$ws.Range("J43").Value = 7432
$ws.Range("L43").Value = 7432
$ws.Range("N43").Value = -7570

$ws.Range("H64").Value = 9352.733
$ws.Range("J64").Value = 10699.8
$ws.Range("L64").Value = 10699.8
$ws.Range("N64").Value = -11195.8

$ws.Range("H67").Value = 9352.733
$ws.Range("J67").Value = 10699.8
$ws.Range("L67").Value = 10699.8
$ws.Range("N67").Value = -12415.8

$ws.Range("H137").Value = 18521216
$ws.Range("I137").Value = 43480468
$ws.Range("J137").Value = 3062.2258
$ws.Range("K137").Value = 130441404
$ws.Range("L137").Value = 9186.6774
$ws.Range("M137").Value = -130438854
$ws.Range("N137").Value = -14286.6774

$ws.Range("H138").Value = 4153.5938
$ws.Range("I138").Value = 3435.5
$ws.Range("K138").Value = 10306.5
$ws.Range("M138").Value = -5166.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 40006136
$ws.Range("I32").Value = 45459336
$ws.Range("J32").Value = 16004.667
$ws.Range("K32").Value = 45459336
$ws.Range("L32").Value = 16004.667
$ws.Range("M32").Value = -45459049
$ws.Range("N32").Value = -16578.667

$ws.Range("H43").Value = 41245.5
$ws.Range("I43").Value = 33833
$ws.Range("K43").Value = 33833
$ws.Range("M43").Value = -33520

$ws.Range("H61").Value = 8216.107
$ws.Range("I61").Value = 5712.675
$ws.Range("K61").Value = 5712.675
$ws.Range("M61").Value = -5500.675

$ws.Range("H136").Value = 8216.107
$ws.Range("I136").Value = 5712.675
$ws.Range("K136").Value = 17138.025
$ws.Range("M136").Value = -14588.025

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5443.1665
$ws.Range("I105").Value = 2411.85
$ws.Range("K105").Value = 2411.85
$ws.Range("M105").Value = -664.8499999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 40499.5
$ws.Range("J28").Value = 40499.5
$ws.Range("L28").Value = 40499.5
$ws.Range("N28").Value = -40989.5

$ws.Range("H76").Value = 4999.8
$ws.Range("I76").Value = 4999.8
$ws.Range("K76").Value = 4999.8
$ws.Range("M76").Value = -4684.8

$ws.Range("H79").Value = 4999.8
$ws.Range("I79").Value = 4999.8
$ws.Range("K79").Value = 4999.8
$ws.Range("M79").Value = -3907.8

$ws.Range("H107").Value = 1207.0834
$ws.Range("I107").Value = 1407.5
$ws.Range("J107").Value = 1006.6667
$ws.Range("K107").Value = 1407.5
$ws.Range("L107").Value = 1006.6667
$ws.Range("M107").Value = 512.5
$ws.Range("N107").Value = -4846.6667

$ws.Range("H141").Value = 144994.5
$ws.Range("J141").Value = 144994.5
$ws.Range("L141").Value = 144994.5
$ws.Range("N141").Value = -155354.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 1080
$ws.Range("I16").Value = 200
$ws.Range("J16").Value = 1666.6666
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 4999.9998
$ws.Range("M16").Value = -427
$ws.Range("N16").Value = -5345.9998

$ws.Range("H17").Value = 379.33334
$ws.Range("I17").Value = 100
$ws.Range("J17").Value = 519
$ws.Range("K17").Value = 300
$ws.Range("L17").Value = 1557
$ws.Range("M17").Value = -131
$ws.Range("N17").Value = -1895

$ws.Range("H43").Value = 1158
$ws.Range("I43").Value = 200
$ws.Range("J43").Value = 4990
$ws.Range("K43").Value = 600
$ws.Range("L43").Value = 14970
$ws.Range("M43").Value = -486
$ws.Range("N43").Value = -15198

$ws.Range("H81").Value = 3929.4285
$ws.Range("J81").Value = 5161.923
$ws.Range("L81").Value = 15485.769
$ws.Range("N81").Value = -17731.769

$ws.Range("H84").Value = 3929.4285
$ws.Range("J84").Value = 5161.923
$ws.Range("L84").Value = 46457.307
$ws.Range("N84").Value = -57689.307

$ws.Range("H86").Value = 643.4737
$ws.Range("J86").Value = 1118.8889
$ws.Range("L86").Value = 3356.6667
$ws.Range("N86").Value = -5728.6667

$ws.Range("H89").Value = 643.4737
$ws.Range("J89").Value = 1118.8889
$ws.Range("L89").Value = 10070.0001
$ws.Range("N89").Value = -21926.0001

$ws.Range("H122").Value = 8870344
$ws.Range("I122").Value = 6410904.5
$ws.Range("J122").Value = 9526195
$ws.Range("K122").Value = 57698140.5
$ws.Range("L122").Value = 85735755
$ws.Range("M122").Value = -57695690.5
$ws.Range("N122").Value = -85740655

$ws.Range("H131").Value = 6656908.5
$ws.Range("J131").Value = 5257144
$ws.Range("L131").Value = 15771432
$ws.Range("N131").Value = -15781512

$ws.Range("H140").Value = 3122.8333
$ws.Range("I140").Value = 1247.2307
$ws.Range("K140").Value = 3741.6921
$ws.Range("M140").Value = 1438.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125").Value = 49442
$ws.Range("J125").Value = 49442
$ws.Range("L125").Value = 49442
$ws.Range("N125").Value = -54362

$ws.Range("H132").Value = 271083.53
$ws.Range("I132").Value = 322941.16
$ws.Range("K132").Value = 968823.48
$ws.Range("M132").Value = -966293.48

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 500001340
$ws.Range("I13").Value = 500001340
$ws.Range("K13").Value = 500001340
$ws.Range("M13").Value = -500001200

$ws.Range("H43").Value = 14671.333
$ws.Range("J43").Value = 24014
$ws.Range("L43").Value = 24014
$ws.Range("N43").Value = -24400

$ws.Range("H46").Value = 4422.3335
$ws.Range("I46").Value = 900
$ws.Range("J46").Value = 5428.7144
$ws.Range("K46").Value = 900
$ws.Range("L46").Value = 5428.7144
$ws.Range("M46").Value = -712
$ws.Range("N46").Value = -5804.7144

$ws.Range("H132").Value = 6702
$ws.Range("I132").Value = 2999.25
$ws.Range("K132").Value = 8997.75
$ws.Range("M132").Value = -6467.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5265.5
$ws.Range("I62").Value = 6143.5
$ws.Range("J62").Value = 4387.5
$ws.Range("K62").Value = 6143.5
$ws.Range("L62").Value = 4387.5
$ws.Range("M62").Value = -5519.5
$ws.Range("N62").Value = -5635.5

$ws.Range("H65").Value = 5265.5
$ws.Range("I65").Value = 6143.5
$ws.Range("J65").Value = 4387.5
$ws.Range("K65").Value = 30717.5
$ws.Range("L65").Value = 21937.5
$ws.Range("M65").Value = -27597.5
$ws.Range("N65").Value = -28177.5

$ws.Range("H131").Value = 113230
$ws.Range("J131").Value = 113230
$ws.Range("L131").Value = 113230
$ws.Range("N131").Value = -123310

$ws.Range("H132").Value = 4166.052
$ws.Range("I132").Value = 1800.6595
$ws.Range("J132").Value = 14272.728
$ws.Range("K132").Value = 5401.9785
$ws.Range("L132").Value = 42818.18399999999
$ws.Range("M132").Value = -2871.9785
$ws.Range("N132").Value = -47878.18399999999
